$wb = $excel.ActiveWorkbook

# --- Overview sheet: shares the "Status" shared string with zh-cn/de-de,
# so it must be updated too for the shared string table to collapse back
# to a single entry (matches upstream, which only touches sharedStrings.xml).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")

# Status column (C2/C3): "Ready for handoff" -> "Handed back: in sync with en-US"
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

# New "Latest Target File" (F) / "Latest Handback File" (G) columns for row 2
$zhcn.Range("F2").Value = "e03118b0-f6e6-4ce8-9ae8-b737145ccbda.md"
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/49f8a99768b8237fe4db933080c5c7caa515fa9d/e2e/e03118b0-f6e6-4ce8-9ae8-b737145ccbda.md", "", "", "e03118b0-f6e6-4ce8-9ae8-b737145ccbda.md") | Out-Null

$zhcn.Range("G2").Value = "e03118b0-f6e6-4ce8-9ae8-b737145ccbda.2f99f4e33ba8012250f4f99a76931040f084d52e.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1068b522949411f1a2eea8dd47fa10e012ab51c7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e03118b0-f6e6-4ce8-9ae8-b737145ccbda.2f99f4e33ba8012250f4f99a76931040f084d52e.zh-cn.xlf", "", "", "e03118b0-f6e6-4ce8-9ae8-b737145ccbda.2f99f4e33ba8012250f4f99a76931040f084d52e.zh-cn.xlf") | Out-Null

# Latest Handback DateTime (H2)
$zhcn.Range("H2").Value = "2016-03-24 07:14:14"

# Row 3
$zhcn.Range("F3").Value = "f4bed4e5-0eb7-4f37-a775-b91429a4a7c2.md"
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/49f8a99768b8237fe4db933080c5c7caa515fa9d/e2e/f4bed4e5-0eb7-4f37-a775-b91429a4a7c2.md", "", "", "f4bed4e5-0eb7-4f37-a775-b91429a4a7c2.md") | Out-Null

$zhcn.Range("G3").Value = "f4bed4e5-0eb7-4f37-a775-b91429a4a7c2.30839b9e218fc24b43be16d4a8ead2d761d7fe8d.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1068b522949411f1a2eea8dd47fa10e012ab51c7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f4bed4e5-0eb7-4f37-a775-b91429a4a7c2.30839b9e218fc24b43be16d4a8ead2d761d7fe8d.zh-cn.xlf", "", "", "f4bed4e5-0eb7-4f37-a775-b91429a4a7c2.30839b9e218fc24b43be16d4a8ead2d761d7fe8d.zh-cn.xlf") | Out-Null

# Latest Handback DateTime (H3)
$zhcn.Range("H3").Value = "2016-03-24 07:14:14"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")

# Status column (C2/C3)
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

# New "Latest Target File" (F) / "Latest Handback File" (G) columns for row 2
$dede.Range("F2").Value = "e03118b0-f6e6-4ce8-9ae8-b737145ccbda.md"
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/49f8a99768b8237fe4db933080c5c7caa515fa9d/e2e/e03118b0-f6e6-4ce8-9ae8-b737145ccbda.md", "", "", "e03118b0-f6e6-4ce8-9ae8-b737145ccbda.md") | Out-Null

$dede.Range("G2").Value = "e03118b0-f6e6-4ce8-9ae8-b737145ccbda.2f99f4e33ba8012250f4f99a76931040f084d52e.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7b23afbfb95441064f4bcc862ebe5f6039f0d7e8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e03118b0-f6e6-4ce8-9ae8-b737145ccbda.2f99f4e33ba8012250f4f99a76931040f084d52e.de-de.xlf", "", "", "e03118b0-f6e6-4ce8-9ae8-b737145ccbda.2f99f4e33ba8012250f4f99a76931040f084d52e.de-de.xlf") | Out-Null

# Latest Handback DateTime (H2)
$dede.Range("H2").Value = "2016-03-24 07:14:21"

# Row 3
$dede.Range("F3").Value = "f4bed4e5-0eb7-4f37-a775-b91429a4a7c2.md"
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/49f8a99768b8237fe4db933080c5c7caa515fa9d/e2e/f4bed4e5-0eb7-4f37-a775-b91429a4a7c2.md", "", "", "f4bed4e5-0eb7-4f37-a775-b91429a4a7c2.md") | Out-Null

$dede.Range("G3").Value = "f4bed4e5-0eb7-4f37-a775-b91429a4a7c2.30839b9e218fc24b43be16d4a8ead2d761d7fe8d.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7b23afbfb95441064f4bcc862ebe5f6039f0d7e8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f4bed4e5-0eb7-4f37-a775-b91429a4a7c2.30839b9e218fc24b43be16d4a8ead2d761d7fe8d.de-de.xlf", "", "", "f4bed4e5-0eb7-4f37-a775-b91429a4a7c2.30839b9e218fc24b43be16d4a8ead2d761d7fe8d.de-de.xlf") | Out-Null

# Latest Handback DateTime (H3)
$dede.Range("H3").Value = "2016-03-24 07:14:21"
